# Update "paises.xlsx" covid stats + re-sort by total cases (column B) descending.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last refreshed" timestamp shown in A1.
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 11:29"

# 2) Update the per-country rows that received new numbers in this refresh.
#    (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$updates = @(
    @{ Pais = "Iran";                 B = 35408; C = 3076; D = 11679; E = 21212; F = 3206; G = 139; H = 2517 },
    @{ Pais = "Suiza";                B = 13187; C = 259;  D = 1530;  E = 11417; F = 203;  G = 9;   H = 240  },
    @{ Pais = "Belgica";              B = 9134;  C = 1850; D = 1063;  E = 7718;  F = 789;  G = 64;  H = 353  },
    @{ Pais = "Noruega";              B = 3796;  C = 25;   D = 7;     E = 3769;  F = 76;   G = 1;   H = 20   },
    @{ Pais = "Eslovenia";            B = 632;   C = 0;    D = 10;    E = 613;   F = 25;   G = 0;   H = 9    },
    @{ Pais = "Sri Lanka";            B = 109;   C = 3;    D = 9;     E = 100;   F = 5;    G = 0;   H = 0    },
    @{ Pais = "Estado de Palestina";  B = 97;    C = 6;    D = 18;    E = 78;    F = 0;    G = 0;   H = 1    }
)

foreach ($u in $updates) {
    $cell = $ws.Range("A4:A205").Find($u.Pais)
    $r = $cell.Row
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
    $ws.Cells.Item($r, 6).Value = $u.F
    $ws.Cells.Item($r, 7).Value = $u.G
    $ws.Cells.Item($r, 8).Value = $u.H
}

# 3) Re-sort the country table (rows 4-205) by "Casos totales" (column B) descending,
#    same as the source dashboard re-ranking countries after each refresh.
$sortRange = $ws.Range("A4:H205")
$keyRange = $ws.Range("B4:B205")
$sortRange.Sort($keyRange, 2)
